$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# Note: the host's ColumnWidth setter round-trips through an internal pixel
# representation (quantized in 1/6-character steps), so the inputs below
# are chosen to land as closely as possible on the target OOXML widths
# (14.4, 6, 14.4, 10.8, 10.8, 9.6, 10.8, 15.6).
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 5.166666666666667
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 10
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 8.833333333333334
$ws.Columns.Item(7).ColumnWidth = 10
$ws.Columns.Item(8).ColumnWidth = 14.833333333333334

# --- Header row (row 1): replace spaces with line breaks ---
$ws.Range("B1").Value = "Kia`nRio"
$ws.Range("C1").Value = "Volkswagen`nGolf"
$ws.Range("D1").Value = "Toyota`nCorolla"
$ws.Range("E1").Value = "Skoda`nOctavia"
$ws.Range("F1").Value = "BMW`n3`nSeries"
$ws.Range("G1").Value = "Hyundai`nSolaris"
$ws.Range("H1").Value = "Вектор`nприоритетов"

# --- Row labels (column A): replace spaces with line breaks ---
$ws.Range("A2").Value = "Kia`nRio"
$ws.Range("A3").Value = "Volkswagen`nGolf"
$ws.Range("A4").Value = "Toyota`nCorolla"
$ws.Range("A5").Value = "Skoda`nOctavia"
$ws.Range("A6").Value = "BMW`n3`nSeries"
$ws.Range("A7").Value = "Hyundai`nSolaris"

# AutoFit the rows that now contain wrapped multi-line text so the engine
# doesn't stamp an explicit customHeight on them (matches the diff, which
# leaves row heights on their default).
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()

# --- Updated priority vector values ---
$ws.Range("H2").Value = "0.276"
$ws.Range("H3").Value = "0.092"
$ws.Range("H4").Value = "0.157"
$ws.Range("H5").Value = "0.157"
$ws.Range("H6").Value = "0.043"
$ws.Range("H7").Value = "0.276"
$ws.Range("H8").Value = "λ_max = 6.062"
$ws.Range("H9").Value = "ИС = 0.012"
